$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.101938724517822
$ws.Range("B1").Value = 2.245247602462769
$ws.Range("C1").Value = 9.548543930053711
$ws.Range("D1").Value = 2.240005970001221
$ws.Range("E1").Value = 1.28649640083313
